# Mini project #1 completed: update the "Status" column (E) on the "Req"
# sheet from "Not implemented" to either "Released" (top-level requirement
# headers) or "Implemented" (sub-requirement detail rows). Rows that were
# already "Implemented" stay "Implemented".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Req")
$ws.Activate()

$statusByRow = @{
    2  = "Released"
    3  = "Released"
    4  = "Released"
    5  = "Implemented"
    6  = "Implemented"
    7  = "Implemented"
    8  = "Released"
    9  = "Implemented"
    10 = "Released"
    11 = "Implemented"
    12 = "Implemented"
    13 = "Implemented"
    14 = "Released"
    15 = "Implemented"
    16 = "Released"
    17 = "Implemented"
    18 = "Implemented"
    19 = "Released"
    20 = "Implemented"
    21 = "Implemented"
    22 = "Released"
    23 = "Implemented"
    24 = "Implemented"
    25 = "Implemented"
    26 = "Released"
    27 = "Implemented"
    28 = "Implemented"
    29 = "Released"
    30 = "Released"
    31 = "Released"
    32 = "Released"
    33 = "Released"
    34 = "Released"
    35 = "Released"
    36 = "Released"
    37 = "Released"
    38 = "Released"
    39 = "Released"
    40 = "Released"
    41 = "Released"
    42 = "Released"
    43 = "Released"
    44 = "Released"
    45 = "Released"
    46 = "Released"
    47 = "Released"
    48 = "Released"
    49 = "Released"
    50 = "Released"
    51 = "Released"
    52 = "Released"
}

foreach ($row in $statusByRow.Keys) {
    $ws.Range("E$row").Value = $statusByRow[$row]
}

# View-state bookkeeping to mirror the saved scroll/selection position.
$ws.Range("F36").Select()
